# Skill.xlsx: "unify the conception of DataNode, DataTable, Entity"
#
# The only functional edits in the target diff (the rest of the diff is
# boilerplate that a different Excel build/locale stamps on every resave -
# fileVersion/rupBuild, xr*/x15 GUID bookkeeping, absPath, window chrome,
# float-precision column-width churn from a Mac->Windows re-save, the
# localized "Normal"->"常规" cell-style name, and the auto phoneticPr/font
# Excel adds once it notices CJK shared strings) are:
#
#   1. The worksheet "Property1" is renamed to "DataNode".
#   2. The active cell in the frozen (bottom) pane moves from K9 to O40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename sheet "Property1" -> "DataNode"
$ws.Name = "DataNode"

# 2) Move the selection in the bottom-left (frozen) pane to O40
$ws.Range("O40").Select() | Out-Null
